$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.577.15"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "2.614.84"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'536.30"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'142.94"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +3.43%  "
$ws.Range("D9").Value = "2.619.11"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  +4.08%  "
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("D14").Value = "3.073.77"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "58.524.00"
$ws.Range("D16").Value = "'20.81"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "2.602.21"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "'335.16"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "'10.15"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").Value = "'6.22"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'67.04"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  +3.26%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").Value = "0.0₃0737"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "'5.92"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").Value = "'153.05"
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").Value = "'0.837"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").Value = "'3.61"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("D42").Value = "'284.07"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'0.594"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("D45").Value = "'10.69"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").Value = "'19.12"
$ws.Range("E47").Value = "  +2.41%  "
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").Value = "'0.0227"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("D50").Value = "1.941.62"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "'4.48"
$ws.Range("E51").Value = "  -0.56%  "
